# New crime data collected (CompStat weekly report roll-forward)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text: volume/issue number and report week dates (rich-text shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  8"
$ws.Range("C9").Value = "Report Covering the Week  2/17/2025  Through  2/23/2025"

# --- Column E best-fit width grows because new % values need an extra digit ---
$ws.Columns.Item(5).ColumnWidth = 7.433768

# --- Row 15 (Rape) ---
$ws.Range("L15").Value = -25
$ws.Range("N15").Value = -76.923076923076

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 29
$ws.Range("K16").Value = 17.241379310344
$ws.Range("L16").Value = 36
$ws.Range("M16").Value = -2.857142857142
$ws.Range("N16").Value = -68.518518518518

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 17.391304347826
$ws.Range("I17").Value = 55
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = 3.773584905660
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 175
$ws.Range("N17").Value = -28.571428571428

# --- Row 18 (Burglary) ---
$ws.Range("C18").Value = 1
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -75
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = -35
$ws.Range("L18").Value = -31.578947368421
$ws.Range("M18").Value = 8.333333333333
$ws.Range("N18").Value = -79.365079365079

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 27
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 22.727272727272
$ws.Range("I19").Value = 54
$ws.Range("J19").Value = 54
$ws.Range("L19").Value = -21.739130434782
$ws.Range("M19").Value = 68.75
$ws.Range("N19").Value = -44.329896907216

# --- Row 20 (G.L.A.) ---
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -40
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 7
$ws.Range("K20").Value = -42.857142857142
$ws.Range("L20").Value = -73.333333333333
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -92.982456140350

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = 8.108108108108
$ws.Range("I21").Value = 164
$ws.Range("J21").Value = 165
$ws.Range("K21").Value = -0.606060606060
$ws.Range("L21").Value = -11.351351351351
$ws.Range("M21").Value = 57.692307692307
$ws.Range("N21").Value = -60.952380952380

# --- Row 22 (Transit): F22 flips from a numeric 1 to the text "0" (no activity) ---
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "0"
$ws.Range("G22").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -40
$ws.Range("F23").Value = 26
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 13.043478260869
$ws.Range("I23").Value = 47
$ws.Range("J23").Value = 51
$ws.Range("K23").Value = -7.843137254901
$ws.Range("L23").Value = -7.843137254901
$ws.Range("M23").Value = 74.074074074074

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 28
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 65.079365079365
$ws.Range("I24").Value = 156
$ws.Range("J24").Value = 108
$ws.Range("K24").Value = 44.444444444444
$ws.Range("L24").Value = 39.285714285714
$ws.Range("M24").Value = 75.280898876404

# --- Row 25 (Retail Theft) ---
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 900
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 825
$ws.Range("I25").Value = 48
$ws.Range("J25").Value = 14
$ws.Range("K25").Value = 242.857142857143
$ws.Range("L25").Value = 200

# --- Row 26 (Misd. Assault) ---
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 18
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = -23.214285714285
$ws.Range("I26").Value = 77
$ws.Range("J26").Value = 87
$ws.Range("K26").Value = -11.494252873563
$ws.Range("L26").Value = -4.938271604938
$ws.Range("M26").Value = -6.097560975609

# --- Row 27 (UCR Rape*) ---
$ws.Range("L27").Value = -55.555555555555

# --- Row 28 (Other Sex Crimes): D28/E28 flip from "n/a" text to real numbers ---
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$excel.CutCopyMode = $false
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 7
$ws.Range("J28").Value = 8
$ws.Range("K28").Value = -12.5
$ws.Range("L28").Value = -30
